$d = $word.ActiveDocument

# Locate the two consecutive "Personal Profile" body paragraphs that are
# being merged/rewritten:
#   P1: "My pursuit of psychology at the degree level ... Leveraging my psychological"
#   P2: "background, my aim is to approach software development ... user needs."
$p1 = $null
$p2 = $null
foreach ($p in $d.Paragraphs) {
    if ($p.Range.Text.StartsWith("My pursuit of psychology at the degree level")) {
        $p1 = $p
    }
    if ($p.Range.Text.StartsWith("background, my aim is to approach software")) {
        $p2 = $p
    }
}

if ($p1 -eq $null -or $p2 -eq $null) {
    throw "Could not locate the Personal Profile paragraphs to replace"
}

# Range spanning from the start of the first paragraph through the end of
# the second paragraph's text (but not its paragraph mark), so the two
# paragraphs collapse into a single merged paragraph.
$start = $p1.Range.Start
$end = $p2.Range.End - 1
$target = $d.Range($start, $end)

$rPr = '<w:rPr><w:rFonts w:ascii="Calibri" w:hAnsi="Calibri" w:eastAsia="Calibri" w:cs="Calibri"/><w:b w:val="0"/><w:bCs w:val="0"/><w:i w:val="0"/><w:iCs w:val="0"/><w:caps w:val="0"/><w:smallCaps w:val="0"/><w:noProof w:val="0"/><w:color w:val="000000" w:themeColor="text1" w:themeTint="FF" w:themeShade="FF"/><w:sz w:val="22"/><w:szCs w:val="22"/><w:lang w:val="en-GB"/></w:rPr>'

$run1 = '<w:r>' + $rPr + '<w:t xml:space="preserve">I pursued psychology at university driven by a curiosity to understand myself and human </w:t></w:r>'
$run2 = '<w:r>' + $rPr + '<w:t>behaviour</w:t></w:r>'
$run3 = '<w:r>' + $rPr + '<w:t xml:space="preserve">, finding fulfilment in unravelling the complexities of social dynamics, body language, and communication processes. My fascination lies in the interconnectedness of underlying processes and observable experiences, reflecting a lifelong interest in the </w:t></w:r>'
$run4 = '<w:r>' + $rPr + '<w:t>synergy</w:t></w:r>'
$run5 = '<w:r>' + $rPr + '<w:t xml:space="preserve"> between back-end mechanisms and front-end interactions.</w:t></w:r>'

$pPr = '<w:pPr><w:pStyle w:val="Normal"/><w:spacing w:after="0" w:line="240" w:lineRule="auto"/><w:rPr/></w:pPr>'

$paragraphXml = '<w:p>' + $pPr + $run1 + $run2 + $run3 + $run4 + $run5 + '</w:p>'

$packageXml = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' + `
    '<w:body>' + $paragraphXml + '</w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'

$target.InsertXML($packageXml)

Write-Output "Merged paragraph text:"
Write-Output $p1.Range.Text
